$wb = $excel.ActiveWorkbook

$changes = @{
    "展览" = @{
        4  = 398
        5  = 4995
        6  = 4995
        7  = 58
        9  = 492
        10 = 1129
        11 = 685
        12 = 4786
        13 = 18
        14 = 36
        15 = 61
        16 = 196
        17 = 202
        18 = 91
        20 = 3695
        23 = 33
        24 = 3501
        25 = 161
        26 = 153
        28 = 179
        29 = 225
        31 = 101
        32 = 95
        33 = 37
        36 = 6142
        37 = 970
        38 = 466
        39 = 92
        40 = 965
        42 = 1273
        43 = 141
        44 = 602
        45 = 22
        46 = 2135
        47 = 310
        48 = 84
        49 = 745
        50 = 889
    }
    "演出" = @{
        7  = 129
        8  = 48
        23 = 790
    }
    "全部类型" = @{
        5  = 398
        6  = 4995
        7  = 4996
        8  = 58
        9  = 48
        12 = 492
        13 = 1129
        14 = 685
        15 = 4787
        16 = 18
        17 = 36
        18 = 61
        19 = 196
        20 = 203
        21 = 91
        23 = 3695
        24 = 3501
        25 = 161
        26 = 153
        27 = 179
        28 = 225
        30 = 101
        31 = 95
        35 = 6142
        36 = 970
        37 = 466
        40 = 92
        41 = 965
        42 = 1273
        43 = 141
        44 = 602
        45 = 2135
        46 = 310
        48 = 745
        49 = 889
    }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $changes[$sheetName]
    foreach ($row in $rowMap.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowMap[$row]
    }
}
